# Generate Report for Handback
# Updates the Handoff/Handback datetimes for the "328bd..." row on the
# zh-cn and de-de sheets, reflecting a fresh report generation run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-02-22 17:37:43"
$wsZhCn.Range("G2").Value = "2016-02-22 17:38:28"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-02-22 17:37:54"
$wsDeDe.Range("G2").Value = "2016-02-22 17:38:47"
